$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Moorings")
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

# Order matters for shared-string table layout: write values in the same
# order the original author entered them (E11, A2, E9, E7, E2, then the
# remaining duplicate "Mooring OOIBARCODE" / "Sensor OOIBARCODE" cells).

$ws2.Range("E11").Value = "OL000128"
$ws2.Range("E11").Style = "Normal"

$ws1.Range("A2").Value = "A00380"
$ws1.Range("A2").Style = "Normal"

$ws2.Range("E9").Value = "N00635"
$ws2.Range("E9").Style = "Normal"

$ws2.Range("E7").Value = "N00639"
$ws2.Range("E7").Style = "Normal"

$ws2.Range("E2").Value = "N00637"
$ws2.Range("E2").Style = "Normal"

$ws2.Range("B2").Value = "A00380"
$ws2.Range("B2").Style = "Normal"

$ws2.Range("B3").Value = "A00380"
$ws2.Range("B3").Style = "Normal"
$ws2.Range("E3").Value = "N00637"
$ws2.Range("E3").Style = "Normal"

$ws2.Range("B4").Value = "A00380"
$ws2.Range("B4").Style = "Normal"
$ws2.Range("E4").Value = "N00637"
$ws2.Range("E4").Style = "Normal"

$ws2.Range("B5").Value = "A00380"
$ws2.Range("B5").Style = "Normal"
$ws2.Range("E5").Value = "N00637"
$ws2.Range("E5").Style = "Normal"

$ws2.Range("B7").Value = "A00380"
$ws2.Range("B7").Style = "Normal"

$ws2.Range("B9").Value = "A00380"
$ws2.Range("B9").Style = "Normal"

$ws2.Range("B11").Value = "A00380"
$ws2.Range("B11").Style = "Normal"

# Restore per-sheet selections, then leave Asset_Cal_Info as the active tab.
[void]$ws1.Range("C32").Select()
[void]$ws2.Range("E25").Select()
[void]$ws2.Activate()
